# Insert a new data row at row 457 (pushing the existing rows 457..558
# down to 458..559), mirroring how this row was added to the underlying
# "Vega Modelo de Temuco - Zanahoria" daily price dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 457; Excel shifts rows 457..558
# down to 458..559 automatically (xlShiftDown is the default behavior).
$ws.Rows(457).Insert()

# Populate the newly inserted row 457 with the new record's data.
$ws.Range("A457").Value = 10
$ws.Range("B457").Value = "Vega Modelo de Temuco"
$ws.Range("C457").Value = "La Araucanía"
$ws.Range("D457").Value = 45244
$ws.Range("E457").Value = 9
$ws.Range("F457").Value = 100114013
$ws.Range("G457").Value = "Zanahoria"
$ws.Range("H457").Value = "Sin especificar"
$ws.Range("I457").Value = "Primera"
$ws.Range("J457").Value = 115
$ws.Range("K457").Value = 6000
$ws.Range("L457").Value = 6000
$ws.Range("M457").Value = 6000
$ws.Range("N457").Value = "$/saco 20 kilos"
$ws.Range("O457").Value = "Región del Maule"
$ws.Range("P457").Value = 300
$ws.Range("Q457").Value = 20
$ws.Range("R457").Value = "Hortaliza"
